# Applies updated Leve profit-calculation figures (columns H-N) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets, per the scheduled-runner
# price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4684.769
$ws.Range("I19").Value = 6690
$ws.Range("J19").Value = 2345.3333
$ws.Range("K19").Value = 6690
$ws.Range("L19").Value = 2345.3333
$ws.Range("M19").Value = -6515
$ws.Range("N19").Value = -2695.3333

$ws.Range("H62").Value = 2720.2
$ws.Range("I62").Value = 3035
$ws.Range("J62").Value = 2248
$ws.Range("K62").Value = 3035
$ws.Range("L62").Value = 2248
$ws.Range("M62").Value = -2411
$ws.Range("N62").Value = -3496

$ws.Range("H65").Value = 2720.2
$ws.Range("I65").Value = 3035
$ws.Range("J65").Value = 2248
$ws.Range("K65").Value = 15175
$ws.Range("L65").Value = 11240
$ws.Range("M65").Value = -12055
$ws.Range("N65").Value = -17480

$ws.Range("H70").Value = 2483.1667
$ws.Range("I70").Value = 2400
$ws.Range("J70").Value = 2524.75
$ws.Range("K70").Value = 7200
$ws.Range("L70").Value = 7574.25
$ws.Range("M70").Value = -6930
$ws.Range("N70").Value = -8114.25

$ws.Range("H73").Value = 2483.1667
$ws.Range("I73").Value = 2400
$ws.Range("J73").Value = 2524.75
$ws.Range("K73").Value = 7200
$ws.Range("L73").Value = 7574.25
$ws.Range("M73").Value = -6264
$ws.Range("N73").Value = -9446.25

$ws.Range("H76").Value = 252750.75
$ws.Range("I76").Value = 252750.75
$ws.Range("K76").Value = 252750.75
$ws.Range("M76").Value = -252435.75

$ws.Range("H79").Value = 252750.75
$ws.Range("I79").Value = 252750.75
$ws.Range("K79").Value = 252750.75
$ws.Range("M79").Value = -251658.75

$ws.Range("H129").Value = 816360.25
$ws.Range("J129").Value = 1045828.6
$ws.Range("L129").Value = 3137485.8
$ws.Range("N129").Value = -3147485.8

$ws.Range("H138").Value = 2762.99
$ws.Range("I138").Value = 1186.9231
$ws.Range("J138").Value = 3316.7432
$ws.Range("K138").Value = 3560.7693
$ws.Range("L138").Value = 9950.229599999999
$ws.Range("M138").Value = 1579.2307
$ws.Range("N138").Value = -20230.2296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1218
$ws.Range("I45").Value = 1096.6666
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 1096.6666
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -719.6666
$ws.Range("N45").Value = -2154

$ws.Range("H63").Value = 837255.9399999999
$ws.Range("I63").Value = 1113674.5
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 1113674.5
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1112988.5
$ws.Range("N63").Value = -9372

$ws.Range("H66").Value = 837255.9399999999
$ws.Range("I66").Value = 1113674.5
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 5568372.5
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -5564940.5
$ws.Range("N66").Value = -46864

$ws.Range("H122").Value = 2682.4
$ws.Range("I122").Value = 2653
$ws.Range("K122").Value = 7959
$ws.Range("M122").Value = -5509

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22227258
$ws.Range("I31").Value = 2825.7
$ws.Range("J31").Value = 40006804
$ws.Range("K31").Value = 2825.7
$ws.Range("L31").Value = 40006804
$ws.Range("M31").Value = -2530.7
$ws.Range("N31").Value = -40007394

$ws.Range("H34").Value = 22227258
$ws.Range("I34").Value = 2825.7
$ws.Range("J34").Value = 40006804
$ws.Range("K34").Value = 2825.7
$ws.Range("L34").Value = 40006804
$ws.Range("M34").Value = -2623.7
$ws.Range("N34").Value = -40007208

$ws.Range("H62").Value = 166668860
$ws.Range("J62").Value = 333333340
$ws.Range("L62").Value = 333333340
$ws.Range("N62").Value = -333334588

$ws.Range("H65").Value = 166668860
$ws.Range("J65").Value = 333333340
$ws.Range("L65").Value = 1666666700
$ws.Range("N65").Value = -1666672940

$ws.Range("H133").Value = 30790.75
$ws.Range("J133").Value = 30790.75
$ws.Range("L133").Value = 30790.75
$ws.Range("N133").Value = -35850.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1100
$ws.Range("I41").Value = 400
$ws.Range("J41").Value = 1800
$ws.Range("K41").Value = 1200
$ws.Range("L41").Value = 5400
$ws.Range("M41").Value = -862
$ws.Range("N41").Value = -6076

$ws.Range("H122").Value = 1069.5714
$ws.Range("J122").Value = 1141.3125
$ws.Range("L122").Value = 10271.8125
$ws.Range("N122").Value = -15171.8125

$ws.Range("H140").Value = 2375.238
$ws.Range("I140").Value = 1340.8334
$ws.Range("J140").Value = 3754.4443
$ws.Range("K140").Value = 4022.5002
$ws.Range("L140").Value = 11263.3329
$ws.Range("M140").Value = 1157.4998
$ws.Range("N140").Value = -21623.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 24600
$ws.Range("I80").Value = 46500
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 46500
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -45502
$ws.Range("N80").Value = -11996

$ws.Range("H83").Value = 24600
$ws.Range("I83").Value = 46500
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 232500
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -227508
$ws.Range("N83").Value = -59984

$ws.Range("H123").Value = 17388.363
$ws.Range("J123").Value = 17388.363
$ws.Range("L123").Value = 17388.363
$ws.Range("N123").Value = -22288.363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8214.857
$ws.Range("I122").Value = 15168
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 45504
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -43054
$ws.Range("N122").Value = -13900

$ws.Range("H133").Value = 17881
$ws.Range("J133").Value = 17881
$ws.Range("L133").Value = 17881
$ws.Range("N133").Value = -22941

$ws.Range("H139").Value = 69150
$ws.Range("J139").Value = 69150
$ws.Range("L139").Value = 69150
$ws.Range("N139").Value = -79430

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H119").Value = 37350
$ws.Range("I119").Value = 10000
$ws.Range("J119").Value = 64700
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 64700
$ws.Range("M119").Value = -5162
$ws.Range("N119").Value = -74376

$ws.Range("H139").Value = 62000
$ws.Range("J139").Value = 62000
$ws.Range("L139").Value = 62000
$ws.Range("N139").Value = -72280
